# Applies the comments from the 2025-06-20 email (12 points logged in the
# bitacora) to the "Form Responses 1" sheet of the projects workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix logo filename typos (drop the stray trailing "2") -----------------
$ws.Range("Q2").Value = "logo_enviome"          # ENVIOME row
$ws.Range("Q3").Value = "logo_acidomic"         # ACIDOMIC row
$ws.Range("Q4").Value = "logo_divergen"         # DIVERGEN row

# CHALLENGE-2 row (row 5): clear placeholder website, fix logo/funding ids
$ws.Range("K5").ClearContents()                 # websiteProject was "none"
$ws.Range("Q5").Value = "logo_CHALLENGE"        # was logo_CHALLENGE2_2
$ws.Range("R5").Value = "CHALLENGE_MICINN"      # was CHALLENGE2_MICINN

# Remaining rows with a placeholder "none" websiteProject -> blank it out
$ws.Range("K7").ClearContents()                 # BlueDNA row
$ws.Range("K8").ClearContents()                 # InGeNi-Caretta row
$ws.Range("K9").ClearContents()                 # GenoMarTur row

# Restore the view/selection state recorded after editing ----------------
$ws.Range("Q5").Select()
